# utmb pv pacing chart
# Fill in the previously-blank "Cimarron" checkpoint row (row 10) on Sheet1
# with its distance/elevation/time-of-day pacing data, then leave the
# selection on I10 (matching the author's last-edited cell).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B10").Value = 1522
$ws.Range("C10").Value = 7
$ws.Range("D10").Formula = "=D9+C10"
$ws.Range("E10").Value = 2715
$ws.Range("F10").Value = 2588
$ws.Range("G10").Value = 0.022916666666666669
$ws.Range("H10").Value = 0.36736111111111108

$ws.Range("I10").Select() | Out-Null
